$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.077808610839221
$ws.Range("D2").Value = 1.082667254768235
$ws.Range("E2").Value = 1.081294811952993
$ws.Range("F2").Value = 1.091784320675822
$ws.Range("I2").Value = 1.048372308005105
$ws.Range("J2").Value = 1.082701223347435
$ws.Range("K2").Value = 1.085334875392519
$ws.Range("L2").Value = 1.083966009003443
$ws.Range("M2").Value = 1.094428429498829
$ws.Range("N2").Value = 1.084238783272155

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.079239888658735
$ws.Range("D3").Value = 1.083958539962137
$ws.Range("E3").Value = 1.08256224496061
$ws.Range("F3").Value = 1.093116892664243
$ws.Range("I3").Value = 1.048665424828717
$ws.Range("J3").Value = 1.083790609522272
$ws.Range("K3").Value = 1.086443693800442
$ws.Range("L3").Value = 1.085050772696539
$ws.Range("M3").Value = 1.095580146410417
$ws.Range("N3").Value = 1.085329716500314

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.080165300350881
$ws.Range("D4").Value = 1.084793631795567
$ws.Range("E4").Value = 1.083381946362196
$ws.Range("F4").Value = 1.093978776662275
$ws.Range("I4").Value = 1.048853329377354
$ws.Range("J4").Value = 1.084494304519896
$ws.Range("K4").Value = 1.087160157864051
$ws.Range("L4").Value = 1.085751710794441
$ws.Range("M4").Value = 1.096324450782184
$ws.Range("N4").Value = 1.08603441082555

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.080554175641311
$ws.Range("D5").Value = 1.085144598814698
$ws.Range("E5").Value = 1.083726453531197
$ws.Range("F5").Value = 1.094341025567214
$ws.Range("I5").Value = 1.048931903744065
$ws.Range("J5").Value = 1.084789851408587
$ws.Range("K5").Value = 1.087461119523895
$ws.Range("L5").Value = 1.086046154579532
$ws.Range("M5").Value = 1.096637136208271
$ws.Range("N5").Value = 1.086330377424722

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.080619459813
$ws.Range("D6").Value = 1.085203521633844
$ws.Range("E6").Value = 1.083784292312477
$ws.Range("F6").Value = 1.094401843748932
$ws.Range("I6").Value = 1.048945072085551
$ws.Range("J6").Value = 1.084839458329159
$ws.Range("K6").Value = 1.087511638342713
$ws.Range("L6").Value = 1.08609557955734
$ws.Range("M6").Value = 1.09668962465556
$ws.Range("N6").Value = 1.086380054792811

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.080170497174901
$ws.Range("D7").Value = 1.084798321844489
$ws.Range("E7").Value = 1.083386550054027
$ws.Range("F7").Value = 1.093983617387694
$ws.Range("I7").Value = 1.048854380943694
$ws.Range("J7").Value = 1.084498254753432
$ws.Range("K7").Value = 1.087164180266765
$ws.Range("L7").Value = 1.085755646068812
$ws.Range("M7").Value = 1.096328629760971
$ws.Range("N7").Value = 1.086038366668871

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.078292469992332
$ws.Range("D8").Value = 1.083103747233382
$ws.Range("E8").Value = 1.081723234367191
$ws.Range("F8").Value = 1.092234750332181
$ws.Range("I8").Value = 1.04847173359041
$ws.Range("J8").Value = 1.083069639213829
$ws.Range("K8").Value = 1.085709818041137
$ws.Range("L8").Value = 1.084332814100625
$ws.Range("M8").Value = 1.094817852923183
$ws.Range("N8").Value = 1.084607722331336

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.07497741203653
$ws.Range("D9").Value = 1.080114036354505
$ws.Range("E9").Value = 1.078788945036857
$ws.Range("F9").Value = 1.089149954992949
$ws.Range("I9").Value = 1.047783921401151
$ws.Range("J9").Value = 1.080542812714707
$ws.Range("K9").Value = 1.083139119345927
$ws.Range("L9").Value = 1.08181797951529
$ws.Range("M9").Value = 1.092148371332956
$ws.Range("N9").Value = 1.082077307448755

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072763217146631
$ws.Range("D10").Value = 1.078118211393094
$ws.Range("E10").Value = 1.076830304252484
$ws.Range("F10").Value = 1.087091130254487
$ws.Range("I10").Value = 1.04731621092377
$ws.Range("J10").Value = 1.078851712977888
$ws.Range("K10").Value = 1.081419781190352
$ws.Range("L10").Value = 1.080136094990032
$ws.Range("M10").Value = 1.090363596513963
$ws.Range("N10").Value = 1.080383806156338

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.071803392967472
$ws.Range("D11").Value = 1.077253309094095
$ws.Range("E11").Value = 1.075981561055491
$ws.Range("F11").Value = 1.086199039925806
$ws.Range("I11").Value = 1.047111497201417
$ws.Range("J11").Value = 1.078117846817259
$ws.Range("K11").Value = 1.080673927720116
$ws.Range("L11").Value = 1.07940651086341
$ws.Range("M11").Value = 1.089589506869037
$ws.Range("N11").Value = 1.079648897821585

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.071446705897882
$ws.Range("D12").Value = 1.076931936652048
$ws.Range("E12").Value = 1.075666199714028
$ws.Range("F12").Value = 1.085867582275201
$ws.Range("I12").Value = 1.047035126719118
$ws.Range("J12").Value = 1.077845010314688
$ws.Range("K12").Value = 1.080396674737233
$ws.Range("L12").Value = 1.07913530893235
$ws.Range("M12").Value = 1.089301780427787
$ws.Range("N12").Value = 1.079375673859885

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.071523224075465
$ws.Range("D13").Value = 1.077000877146916
$ws.Range("E13").Value = 1.075733850402511
$ws.Range("F13").Value = 1.085938685465745
$ws.Range("I13").Value = 1.047051523433334
$ws.Range("J13").Value = 1.077903545867894
$ws.Range("K13").Value = 1.080456156003975
$ws.Range("L13").Value = 1.079193491868662
$ws.Range("M13").Value = 1.089363507622013
$ws.Range("N13").Value = 1.079434292540291

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.071773912516147
$ws.Range("D14").Value = 1.077226746602007
$ws.Range("E14").Value = 1.075955495260556
$ws.Range("F14").Value = 1.086171643509516
$ws.Range("I14").Value = 1.047105191147463
$ws.Range("J14").Value = 1.078095299108391
$ws.Range("K14").Value = 1.080651014205976
$ws.Range("L14").Value = 1.07938409736933
$ws.Range("M14").Value = 1.089565727315741
$ws.Range("N14").Value = 1.079626318092385

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071928347813111
$ws.Range("D15").Value = 1.077365897619941
$ws.Range("E15").Value = 1.076092044547726
$ws.Range("F15").Value = 1.086315163845449
$ws.Range("I15").Value = 1.047138213736573
$ws.Range("J15").Value = 1.078213411892813
$ws.Range("K15").Value = 1.080771044850565
$ws.Range("L15").Value = 1.079501508834089
$ws.Range("M15").Value = 1.089690295566329
$ws.Range("N15").Value = 1.079744598610507

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.072826894613195
$ws.Range("D16").Value = 1.07817559704218
$ws.Range("E16").Value = 1.076886618713146
$ws.Range("F16").Value = 1.087150322200194
$ws.Range("I16").Value = 1.047329750799854
$ws.Range("J16").Value = 1.078900383001748
$ws.Range("K16").Value = 1.081469251836993
$ws.Range("L16").Value = 1.080184486993699
$ws.Range("M16").Value = 1.090414943169234
$ws.Range("N16").Value = 1.080432545297214

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.073390240479981
$ws.Range("D17").Value = 1.07868331033326
$ws.Range("E17").Value = 1.077384860000526
$ws.Range("F17").Value = 1.087674029592154
$ws.Range("I17").Value = 1.047449309115726
$ws.Range("J17").Value = 1.079330868124961
$ws.Range("K17").Value = 1.081906849260725
$ws.Range("L17").Value = 1.080612545702827
$ws.Range("M17").Value = 1.090869152377846
$ws.Range("N17").Value = 1.080863641758679

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.073718728412656
$ws.Range("D18").Value = 1.078979383970627
$ws.Range("E18").Value = 1.077675414374085
$ws.Range("F18").Value = 1.087979441131469
$ws.Range("I18").Value = 1.047518834055718
$ws.Range("J18").Value = 1.079581807874509
$ws.Range("K18").Value = 1.082161960657705
$ws.Range("L18").Value = 1.080862098306002
$ws.Range("M18").Value = 1.091133962639611
$ws.Range("N18").Value = 1.081114937871455

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.07383071715497
$ws.Range("D19").Value = 1.079080326212088
$ws.Range("E19").Value = 1.077774475633501
$ws.Range("F19").Value = 1.088083568853476
$ws.Range("I19").Value = 1.047542504447055
$ws.Range("J19").Value = 1.079667345691012
$ws.Range("K19").Value = 1.082248924825003
$ws.Range("L19").Value = 1.080947167908547
$ws.Range("M19").Value = 1.091224235480886
$ws.Range("N19").Value = 1.081200597161469

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.073329809413261
$ws.Range("D20").Value = 1.078628844459586
$ws.Range("E20").Value = 1.077331409803222
$ws.Range("F20").Value = 1.087617846791803
$ws.Range("I20").Value = 1.047436503516697
$ws.Range("J20").Value = 1.079284697203324
$ws.Range("K20").Value = 1.081859912860222
$ws.Range("L20").Value = 1.080566632200358
$ws.Range("M20").Value = 1.090820432708725
$ws.Range("N20").Value = 1.080817405269038

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.071700095638454
$ws.Range("D21").Value = 1.077160236728475
$ws.Range("E21").Value = 1.075890229178842
$ws.Range("F21").Value = 1.086103045847507
$ws.Range("I21").Value = 1.047089396483001
$ws.Range("J21").Value = 1.078038839381483
$ws.Range("K21").Value = 1.080593639147114
$ws.Range("L21").Value = 1.07932797440608
$ws.Range("M21").Value = 1.089506184075343
$ws.Range("N21").Value = 1.079569778186187

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.070674468657311
$ws.Range("D22").Value = 1.076236230600855
$ws.Range("E22").Value = 1.074983519080122
$ws.Range("F22").Value = 1.085150075068851
$ws.Range("I22").Value = 1.046869242554941
$ws.Range("J22").Value = 1.077254094104786
$ws.Range("K22").Value = 1.079796266718247
$ws.Range("L22").Value = 1.078548010933804
$ws.Range("M22").Value = 1.088678733011708
$ws.Range("N22").Value = 1.078783918481186

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.071218266018404
$ws.Range("D23").Value = 1.076726125355865
$ws.Range("E23").Value = 1.075464239945862
$ws.Range("F23").Value = 1.085655317085282
$ws.Range("I23").Value = 1.046986132197533
$ws.Range("J23").Value = 1.07767023891634
$ws.Range("K23").Value = 1.080219085576959
$ws.Range("L23").Value = 1.078961596596894
$ws.Range("M23").Value = 1.089117489077032
$ws.Range("N23").Value = 1.079200654266105

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.073357115926451
$ws.Range("D24").Value = 1.078653455450039
$ws.Range("E24").Value = 1.077355561835967
$ws.Range("F24").Value = 1.087643233560285
$ws.Range("I24").Value = 1.047442290468827
$ws.Range("J24").Value = 1.079305560333917
$ws.Range("K24").Value = 1.081881121806205
$ws.Range("L24").Value = 1.080587378927382
$ws.Range("M24").Value = 1.090842447407712
$ws.Range("N24").Value = 1.080838298027669

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.075835144401023
$ws.Range("D25").Value = 1.080887404652336
$ws.Range("E25").Value = 1.07954794446747
$ws.Range("F25").Value = 1.089947835894707
$ws.Range("I25").Value = 1.047963348547454
$ws.Range("J25").Value = 1.081197196509149
$ws.Range("K25").Value = 1.08380466743426
$ws.Range("L25").Value = 1.082469048710729
$ws.Range("M25").Value = 1.092839383174863
$ws.Range("N25").Value = 1.082732620543244

